# Update "想去人数" (want-to-go count) values in F column across sheets
# per the gh-pages data regeneration commit 456a3b4.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 1808
$ws.Range("F8").Value = 2414
$ws.Range("F10").Value = 7247
$ws.Range("F11").Value = 176
$ws.Range("F13").Value = 185
$ws.Range("F14").Value = 1652
$ws.Range("F19").Value = 3277
$ws.Range("F20").Value = 5755
$ws.Range("F21").Value = 5755
$ws.Range("F23").Value = 920
$ws.Range("F25").Value = 331
$ws.Range("F26").Value = 5714
$ws.Range("F37").Value = 63
$ws.Range("F40").Value = 466
$ws.Range("F42").Value = 68
$ws.Range("F43").Value = 336

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 3320
$ws.Range("F7").Value = 1512
$ws.Range("F9").Value = 442
$ws.Range("F10").Value = 2972
$ws.Range("F12").Value = 724
$ws.Range("F13").Value = 908
$ws.Range("F14").Value = 906
$ws.Range("F15").Value = 1402

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1512
$ws.Range("F6").Value = 2972
$ws.Range("F7").Value = 1808
$ws.Range("F10").Value = 2414
$ws.Range("F12").Value = 7248
$ws.Range("F13").Value = 176
$ws.Range("F14").Value = 724
$ws.Range("F15").Value = 1652
$ws.Range("F17").Value = 906
$ws.Range("F22").Value = 1402
$ws.Range("F23").Value = 3278
$ws.Range("F25").Value = 5755
$ws.Range("F29").Value = 331
$ws.Range("F30").Value = 5714
$ws.Range("F32").Value = 3935
$ws.Range("F41").Value = 63
$ws.Range("F43").Value = 466
$ws.Range("F45").Value = 336
